# Auto-generated edit script applying numeric updates from the Odin_Profits diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 13006508
$ws.Range("J74").Value = 34374.75
$ws.Range("L74").Value = 34374.75
$ws.Range("N74").Value = -36246.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 13006508
$ws.Range("J77").Value = 34374.75
$ws.Range("L77").Value = 171873.75
$ws.Range("N77").Value = -181233.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1124.4546
$ws.Range("I107").Value = 1124.4546
$ws.Range("K107").Value = 1124.4546
$ws.Range("M107").Value = 795.5454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 307864.16
$ws.Range("I132").Value = 443829.44
$ws.Range("J132").Value = 4557
$ws.Range("K132").Value = 1331488.32
$ws.Range("L132").Value = 13671
$ws.Range("M132").Value = -1328958.32
$ws.Range("N132").Value = -18731

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4552.6553
$ws.Range("I137").Value = 3725.1333
$ws.Range("J137").Value = 5439.2856
$ws.Range("K137").Value = 11175.3999
$ws.Range("L137").Value = 16317.8568
$ws.Range("M137").Value = -8625.3999
$ws.Range("N137").Value = -21417.8568

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5212.8945
$ws.Range("I138").Value = 2072.5
$ws.Range("J138").Value = 5582.353
$ws.Range("K138").Value = 6217.5
$ws.Range("L138").Value = 16747.059
$ws.Range("M138").Value = -1077.5
$ws.Range("N138").Value = -27027.059

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 523.2857
$ws.Range("I141").Value = 523.2857
$ws.Range("K141").Value = 1569.8571
$ws.Range("M141").Value = 3610.1429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3089751
$ws.Range("I32").Value = 2129.6785
$ws.Range("J32").Value = 13896425
$ws.Range("K32").Value = 2129.6785
$ws.Range("L32").Value = 13896425
$ws.Range("M32").Value = -1842.6785
$ws.Range("N32").Value = -13896999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2164.4
$ws.Range("I63").Value = 2011
$ws.Range("J63").Value = 2266.6667
$ws.Range("K63").Value = 2011
$ws.Range("L63").Value = 2266.6667
$ws.Range("M63").Value = -1325
$ws.Range("N63").Value = -3638.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2164.4
$ws.Range("I66").Value = 2011
$ws.Range("J66").Value = 2266.6667
$ws.Range("K66").Value = 10055
$ws.Range("L66").Value = 11333.3335
$ws.Range("M66").Value = -6623
$ws.Range("N66").Value = -18197.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 8948.5
$ws.Range("I88").Value = 2995
$ws.Range("J88").Value = 9610
$ws.Range("K88").Value = 2995
$ws.Range("L88").Value = 9610
$ws.Range("M88").Value = -2589
$ws.Range("N88").Value = -10422

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 8948.5
$ws.Range("I91").Value = 2995
$ws.Range("J91").Value = 9610
$ws.Range("K91").Value = 2995
$ws.Range("L91").Value = 9610
$ws.Range("M91").Value = -1591
$ws.Range("N91").Value = -12418

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 35715884
$ws.Range("I20").Value = 47620344
$ws.Range("K20").Value = 47620344
$ws.Range("M20").Value = -47620097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 34531.332
$ws.Range("I86").Value = 39424.75
$ws.Range("J86").Value = 24744.5
$ws.Range("K86").Value = 39424.75
$ws.Range("L86").Value = 24744.5
$ws.Range("M86").Value = -38301.75
$ws.Range("N86").Value = -26990.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 34531.332
$ws.Range("I89").Value = 39424.75
$ws.Range("J89").Value = 24744.5
$ws.Range("K89").Value = 197123.75
$ws.Range("L89").Value = 123722.5
$ws.Range("M89").Value = -191507.75
$ws.Range("N89").Value = -134954.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5172.0625
$ws.Range("I105").Value = 3887.0908
$ws.Range("K105").Value = 3887.0908
$ws.Range("M105").Value = -2140.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 14566.333
$ws.Range("I62").Value = 17849.5
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 17849.5
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -17225.5
$ws.Range("N62").Value = -9248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 14566.333
$ws.Range("I65").Value = 17849.5
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 89247.5
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -86127.5
$ws.Range("N65").Value = -46240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 45618416
$ws.Range("I131").Value = 53339080
$ws.Range("J131").Value = 37039900
$ws.Range("K131").Value = 160017240
$ws.Range("L131").Value = 111119700
$ws.Range("M131").Value = -160012200
$ws.Range("N131").Value = -111129780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3098.0833
$ws.Range("I132").Value = 1224.75
$ws.Range("J132").Value = 4034.75
$ws.Range("K132").Value = 11022.75
$ws.Range("L132").Value = 36312.75
$ws.Range("M132").Value = -8492.75
$ws.Range("N132").Value = -41372.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 9249.666999999999
$ws.Range("I138").Value = 9249.666999999999
$ws.Range("K138").Value = 27749.001
$ws.Range("M138").Value = -22609.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1944
$ws.Range("I139").Value = 1320.2667
$ws.Range("K139").Value = 3960.800099999999
$ws.Range("M139").Value = 1179.199900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13220.523
$ws.Range("I70").Value = 8881.916999999999
$ws.Range("K70").Value = 8881.916999999999
$ws.Range("M70").Value = -8611.916999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 13220.523
$ws.Range("I73").Value = 8881.916999999999
$ws.Range("K73").Value = 8881.916999999999
$ws.Range("M73").Value = -7945.916999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 27227.8
$ws.Range("I80").Value = 5463
$ws.Range("J80").Value = 59875
$ws.Range("K80").Value = 5463
$ws.Range("L80").Value = 59875
$ws.Range("M80").Value = -4465
$ws.Range("N80").Value = -61871

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 27227.8
$ws.Range("I83").Value = 5463
$ws.Range("J83").Value = 59875
$ws.Range("K83").Value = 27315
$ws.Range("L83").Value = 299375
$ws.Range("M83").Value = -22323
$ws.Range("N83").Value = -309359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2112.0833
$ws.Range("I68").Value = 2312.2222
$ws.Range("J68").Value = 1511.6666
$ws.Range("K68").Value = 2312.2222
$ws.Range("L68").Value = 1511.6666
$ws.Range("M68").Value = -1563.2222
$ws.Range("N68").Value = -3009.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2112.0833
$ws.Range("I71").Value = 2312.2222
$ws.Range("J71").Value = 1511.6666
$ws.Range("K71").Value = 11561.111
$ws.Range("L71").Value = 7558.333000000001
$ws.Range("M71").Value = -7817.111000000001
$ws.Range("N71").Value = -15046.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2876
$ws.Range("J122").Value = 3000
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4492.8486
$ws.Range("I126").Value = 3501.3076
$ws.Range("J126").Value = 8175.7144
$ws.Range("K126").Value = 10503.9228
$ws.Range("L126").Value = 24527.1432
$ws.Range("M126").Value = -8033.9228
$ws.Range("N126").Value = -29467.1432

